$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1 (first table) ---
# Shared-string table order: dw-, dwh, dwf+, w-, wh, wf+, fp, wp, dfp, dwp
# so populate K:M before H:J to match that insertion order.
$ws.Range("K1").Value = "dw-"
$ws.Range("L1").Value = "dwh"
$ws.Range("M1").Value = "dwf+"
$ws.Range("H1").Value = "w-"
$ws.Range("I1").Value = "wh"
$ws.Range("J1").Value = "wf+"
$ws.Range("N1").Value = "fp"
$ws.Range("O1").Value = "wp"
$ws.Range("P1").Value = "dfp"
$ws.Range("Q1").Value = "dwp"

# --- Header row 7 (second table) ---
$ws.Range("K7").Value = "dw-"
$ws.Range("L7").Value = "dwh"
$ws.Range("M7").Value = "dwf+"
$ws.Range("H7").Value = "w-"
$ws.Range("I7").Value = "wh"
$ws.Range("J7").Value = "wf+"

# --- First table: rows 2-5 ---
foreach ($r in 2..5) {
    $ws.Range("H$r").Formula = "=B$r * 2 * PI()"
    $ws.Range("I$r").Formula = "=C$r * 2 * PI()"
    $ws.Range("J$r").Formula = "=D$r * 2 * PI()"
    $ws.Range("K$r").Formula = "=E$r * 2 * PI()"
    $ws.Range("L$r").Formula = "=F$r * 2 * PI()"
    $ws.Range("M$r").Formula = "=G$r * 2 * PI()"
    $ws.Range("N$r").Formula = "=SQRT(D$r*D$r+B$r*B$r-C$r*C$r)"
    $ws.Range("O$r").Formula = "=N$r*2*PI()"
    $ws.Range("P$r").Formula = "=`$N`$2-N$r"
    $ws.Range("Q$r").Formula = "=P$r*2*PI()"
}

# --- Second table: rows 8-11 ---
foreach ($r in 8..11) {
    $ws.Range("H$r").Formula = "=B$r * 2 * PI()"
    $ws.Range("I$r").Formula = "=C$r * 2 * PI()"
    $ws.Range("J$r").Formula = "=D$r * 2 * PI()"
    $ws.Range("K$r").Formula = "=E$r * 2 * PI()"
    $ws.Range("L$r").Formula = "=F$r * 2 * PI()"
    $ws.Range("M$r").Formula = "=G$r * 2 * PI()"
    $ws.Range("N$r").Formula = "=SQRT(D$r*D$r+B$r*B$r-C$r*C$r)"
    $ws.Range("O$r").Formula = "=N$r*2*PI()"
    if ($r -eq 8) {
        $ws.Range("P$r").Formula = "=`$N`$8-N$r"
    } else {
        $ws.Range("P$r").Formula = "=`$N`$2-N$r"
    }
    $ws.Range("Q$r").Formula = "=P$r*2*PI()"
}

# --- Update dimension / selection ---
$ws.Range("S7").Select()
